$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (13) to the mail_template sheet for the "password notification header"
# mail template, reusing the "Password notification" / "body" values already used
# by the existing password_notify row, plus two new strings and a new body text.
$ws.Range("B13").Value() = "password_notify_header"
$ws.Range("C13").Value() = $ws.Range("C12").Value()
$ws.Range("D13").Value() = $ws.Range("D11").Value()
$ws.Range("E13").Value() = "[`${system:site_name}]Password notification header"
$ws.Range("F13").Value() = "***************************************************************************\nThis email is automatically encrypted as an attachment.\nYou will receive a password from the sender later.\n***************************************************************************\n\n"

# Select the newly added row, matching the saved selection state in the file.
$ws.Rows.Item(13).Select()
